$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Update time_taken column (F) values on the data sheet
$ws.Range("F2").Value = "2021-10-05 14:33:23.656024"
$ws.Range("F3").Value = "2021-10-05 14:33:23.656032"
$ws.Range("F4").Value = "2021-10-05 14:33:23.656035"
$ws.Range("F5").Value = "2021-10-05 14:33:23.656037"
$ws.Range("F6").Value = "2021-10-05 14:33:23.656039"
$ws.Range("F7").Value = "2021-10-05 14:33:23.656042"
$ws.Range("F8").Value = "2021-10-05 14:33:23.656044"
$ws.Range("F9").Value = "2021-10-05 14:33:23.656046"
$ws.Range("F10").Value = "2021-10-05 14:33:23.656048"
$ws.Range("F11").Value = "2021-10-05 14:33:23.656050"

# Add metadata sheet after the data sheet
$meta = $wb.Worksheets.Add($null, $ws)
$meta.Name = "metadata"

$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Central Hypoventilation"
$meta.Range("C2").Value = 71
# Force "1.3" to be stored as text (not a number) without leaving a quote-prefix
# style behind: write it with a leading apostrophe, then reset the cell
# formatting back to the plain default (copy from an unstyled cell) so no
# extra style entry lingers.
$meta.Range("D2").Value = "'1.3"
$ws.Range("A1").Copy()
$meta.Range("D2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$meta.Range("E2").Value = "2021-09-22T22:38:35.040330Z"
$meta.Range("F2").Value = "2021-10-05 14:33:23.653399"
$meta.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/71/?format=json"

# Style header row and A2 the same as the data sheet's matching cells (bold,
# bordered, centered/top-aligned header style). Copy/PasteSpecial(Formats)
# reuses the existing style record instead of creating a new one.
$ws.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
